# Weekly refresh of fruit/vegetable prices: the per-week records (rows 2-8)
# get their Fecha/Volumen/Precio columns reshuffled to reflect the latest
# weekly data pull. Columns A,B,C,E,F,G,H,I,N,O,Q,R stay tied to the row;
# only D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows.
$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the current (pre-edit) values for each shuffled column, keyed by row.
$snapshot = @{}
foreach ($row in 2..8) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Destination row -> source row (i.e. destination row now carries the values
# that used to live in the source row).
$rowMap = @{
    2 = 8
    3 = 4
    4 = 6
    5 = 5
    6 = 3
    7 = 2
    8 = 7
}

foreach ($destRow in 2..8) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcValues[$col]
    }
}
